{"js": "// Word JavaScript API (Office.js) script.\n// Body of: async (context) => { ... }\n\n// 1) Insert a new paragraph \"Gymn\u00e1zium\" right before the paragraph\n//    \"P\u016fsobil v\u00a0Lond\u00fdn\u011b jako herec a dramatik\", inheriting the same\n//    paragraph formatting (list style / numbering) via insertParagraph.\nconst body = context.document.body;\nconst hits = body.search(\"P\u016fsobil v\u00a0Lond\u00fdn\u011b jako herec a dramatik\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Target paragraph \"P\u016fsobil v\u00a0Lond\u00fdn\u011b jako herec a dramatik\" not found.');\n}\n\nconst targetRange = hits.items[0];\nconst targetParagraph = targetRange.paragraphs.getFirst();\nconst newParagraph = targetParagraph.insertParagraph(\"Gymn\u00e1zium\", Word.InsertLocation.before);\nawait context.sync();\n\n// 2) Append \", ner\u00fdmovan\u00fd ver\u0161\" as a new run after the existing\n//    \" \u2013 p\u011btistop\u00fd jambick\u00fd ver\u0161\" text, within the \"Blankverse\" paragraph.\nconst verseHits = body.search(\" \u2013 p\u011btistop\u00fd jambick\u00fd ver\u0161\", { matchCase: true });\nverseHits.load(\"items\");\nawait context.sync();\n\nif (verseHits.items.length === 0) {\n  throw new Error('Target text \" \u2013 p\u011btistop\u00fd jambick\u00fd ver\u0161\" not found.');\n}\n\nconst verseRange = verseHits.items[0];\nverseRange.insertText(\", ner\u00fdmovan\u00fd ver\u0161\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Insert a new paragraph \"Gymn\u00e1zium\" right before the paragraph\n#    \"P\u016fsobil v\u00a0Lond\u00fdn\u011b jako herec a dramatik\". InsertParagraphBefore()\n#    splits in a new, empty paragraph right before the found text that\n#    inherits the found paragraph's formatting (list style\n#    \"Odstavecseseznamem\", numbering ilvl 0 / numId 8, spacing after 0);\n#    that new empty paragraph is $rng.Paragraphs(1) after the insert, so we\n#    just set its text.\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"P\u016fsobil v\u00a0Lond\u00fdn\u011b jako herec a dramatik\")\n$rng1.InsertParagraphBefore()\n$newPara = $rng1.Paragraphs(1)\n$newPara.Range.Text = \"Gymn\u00e1zium\"\n\n# 2) Append a new run \", ner\u00fdmovan\u00fd ver\u0161\" right after the existing\n#    \" \u2013 p\u011btistop\u00fd jambick\u00fd ver\u0161\" text in the \"Blankverse\" paragraph.\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\" \u2013 p\u011btistop\u00fd jambick\u00fd ver\u0161\")\n$rng2.InsertAfter(\", ner\u00fdmovan\u00fd ver\u0161\")\n"}
